# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# Sheet "展览" (Exhibition) rows 2, 5, 6 and sheet "全部类型" (All Types) rows 2, 7, 8
# hold the same events, so both need the matching updates.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 385
$ws1.Range("F5").Value = 1050
$ws1.Range("F6").Value = 2439

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 385
$ws4.Range("F7").Value = 1050
$ws4.Range("F8").Value = 2439
